$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.061.08"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "2.416.64"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.351"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.92%  "

$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").Value = "2.854.95"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").Value = "61.995.38"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").Value = "2.414.43"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "321.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.79%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "567.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.97%  "

$ws.Range("D27").Value = "2.533.10"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "0.0₃0933"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("E30").Value = "  -1.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.09%  "

$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.37%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "

$ws.Range("E38").Value = "  -4.24%  "

$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.12%  "

$ws.Range("E41").Value = "  -3.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.89%  "

$ws.Range("E46").Value = "  -2.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.595"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0921"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("E50").Value = "  -0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
